$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update country name cells that shifted rank position (A column) ---
$ws.Range("A34").Value = "Republica Dominicana"
$ws.Range("A35").Value = "Kazajistan"
$ws.Range("A63").Value = "Moldavia"
$ws.Range("A64").Value = "Chequia"
$ws.Range("A65").Value = "Ghana"
$ws.Range("A96").Value = "Namibia"
$ws.Range("A97").Value = "Guinea"
$ws.Range("A122").Value = "Birmania"
$ws.Range("A123").Value = "Republica de Africa Central"
$ws.Range("A124").Value = "Surinam"
$ws.Range("A125").Value = "Ruanda"
$ws.Range("A130").Value = "Trinidad yTobago"
$ws.Range("A131").Value = "Siria"
$ws.Range("A132").Value = "Lituania"
$ws.Range("A204").Value = "Timor Oriental"
$ws.Range("A205").Value = "Santa Lucia"
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("A215").Value = "Montserrat"

# --- Update numeric statistics cells (B,C,D,E,F,G,H columns) ---
$ws.Range("B4").Value = 6932328
$ws.Range("C4").Value = 6387
$ws.Range("D4").Value = 4192963
$ws.Range("E4").Value = 2536091
$ws.Range("G4").Value = 103
$ws.Range("H4").Value = 203274
$ws.Range("B5").Value = 5329983
$ws.Range("C5").Value = 24508
$ws.Range("D5").Value = 4225803
$ws.Range("E5").Value = 1018415
$ws.Range("G5").Value = 140
$ws.Range("H5").Value = 85765
$ws.Range("B20").Value = 315597
$ws.Range("C20").Value = 3907
$ws.Range("D20").Value = 249539
$ws.Range("E20").Value = 57567
$ws.Range("G20").Value = 83
$ws.Range("H20").Value = 8491
$ws.Range("B25").Value = 271672
$ws.Range("C25").Value = 428
$ws.Range("E25").Value = 19206
$ws.Range("B27").Value = 183602
$ws.Range("C27").Value = 4531
$ws.Range("D27").Value = 132449
$ws.Range("E27").Value = 49927
$ws.Range("G27").Value = 30
$ws.Range("H27").Value = 1226
$ws.Range("B34").Value = 107700
$ws.Range("C34").Value = 968
$ws.Range("D34").Value = 80820
$ws.Range("E34").Value = 24836
$ws.Range("G34").Value = 10
$ws.Range("H34").Value = 2044
$ws.Range("B35").Value = 107199
$ws.Range("C35").Value = 65
$ws.Range("D35").Value = 101822
$ws.Range("E35").Value = 3706
$ws.Range("H35").Value = 1671
$ws.Range("B63").Value = 46336
$ws.Range("C63").Value = 688
$ws.Range("D63").Value = 34236
$ws.Range("E63").Value = 10899
$ws.Range("G63").Value = 15
$ws.Range("H63").Value = 1201
$ws.Range("B64").Value = 46262
$ws.Range("C64").Value = 0
$ws.Range("D64").Value = 23858
$ws.Range("E64").Value = 21909
$ws.Range("H64").Value = 495
$ws.Range("B65").Value = 45857
$ws.Range("C65").Value = 97
$ws.Range("D65").Value = 45029
$ws.Range("E65").Value = 533
$ws.Range("H65").Value = 295
$ws.Range("D88").Value = 9989
$ws.Range("E88").Value = 4422
$ws.Range("B94").Value = 12226
$ws.Range("C94").Value = 153
$ws.Range("D94").Value = 6888
$ws.Range("E94").Value = 4980
$ws.Range("G94").Value = 5
$ws.Range("H94").Value = 358
$ws.Range("B96").Value = 10292
$ws.Range("C96").Value = 85
$ws.Range("D96").Value = 7969
$ws.Range("E96").Value = 2212
$ws.Range("G96").Value = 3
$ws.Range("H96").Value = 111
$ws.Range("B97").Value = 10231
$ws.Range("D97").Value = 9660
$ws.Range("E97").Value = 508
$ws.Range("H97").Value = 63
$ws.Range("B122").Value = 4870
$ws.Range("C122").Value = 403
$ws.Range("D122").Value = 1188
$ws.Range("E122").Value = 3601
$ws.Range("G122").Value = 11
$ws.Range("H122").Value = 81
$ws.Range("B123").Value = 4786
$ws.Range("D123").Value = 1830
$ws.Range("E123").Value = 2894
$ws.Range("H123").Value = 62
$ws.Range("B124").Value = 4691
$ws.Range("D124").Value = 4280
$ws.Range("E124").Value = 315
$ws.Range("H124").Value = 96
$ws.Range("B125").Value = 4671
$ws.Range("C125").Value = 0
$ws.Range("D125").Value = 2845
$ws.Range("E125").Value = 1801
$ws.Range("G125").Value = 0
$ws.Range("H125").Value = 25
$ws.Range("B130").Value = 3739
$ws.Range("C130").Value = 88
$ws.Range("D130").Value = 1586
$ws.Range("E130").Value = 2093
$ws.Range("H130").Value = 60
$ws.Range("B131").Value = 3731
$ws.Range("C131").Value = 0
$ws.Range("D131").Value = 918
$ws.Range("E131").Value = 2645
$ws.Range("H131").Value = 168
$ws.Range("B132").Value = 3664
$ws.Range("C132").Value = 99
$ws.Range("D132").Value = 2197
$ws.Range("E132").Value = 1380
$ws.Range("H132").Value = 87
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1

# --- Update "last updated" timestamp message ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Septiembre de 2020 a las 17:01"
